$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting the existing data (and the chart's
# referenced range) down by one row.
$ws.Rows(2).Insert()

# Populate the newly inserted row with the new date entry (2023-03-28, 10 students).
$ws.Range("A2").Value = 45013
$ws.Range("A2").NumberFormat = "yyyy-mm-dd"
$ws.Range("B2").Value = 10

# Update the chart title text.
$chart = $ws.ChartObjects(1).Chart
$chart.ChartTitle.Text = "Students' Entrances"
